$d = $word.ActiveDocument

$d.Content.Find.Execute("Date/Time: 2/2/23 @ 1100", $true, $false, $false, $false, $false, $true, 1, $false, "Date/Time: 2/2/23 @ 1100", 2)
